$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their original Text format so that numeric-looking
# strings (e.g. "305.08") and percentages (e.g. "2.16%") are stored verbatim as text,
# matching the source workbook which stores these columns as inline strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.16%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.12%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.178"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.45%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07517"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.28%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.327"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "32.82%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.021"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.61%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9167"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1735"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07645"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.51%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08283"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.39%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03039"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.61%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09934"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.49%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001518"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.46%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006066"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-6.79%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.496"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.35%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.881"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.32%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.239"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.90%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.86%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.98%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.649"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.76%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04623"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.51%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.02%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.94%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004538"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.19%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002738"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "47.39%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01760"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.01%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04576"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.71%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007228"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.38%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1366"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.82%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002197"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.74%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01078"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-15.66%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006549"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.76%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-57.48%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.009885"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-23.77%"
